$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'311.13"
$ws.Range("E2").Value = "'-0.61%"
$ws.Range("D3").Value = "'37.69"
$ws.Range("E3").Value = "'-0.65%"
$ws.Range("D4").Value = "'5.158"
$ws.Range("E4").Value = "'1.68%"
$ws.Range("D5").Value = "'0.07922"
$ws.Range("E5").Value = "'1.95%"
$ws.Range("B6").Value = "GateToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D6").Value = "'4.434"
$ws.Range("E6").Value = "'1.91%"
$ws.Range("B7").Value = "FTXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D7").Value = "'1.921"
$ws.Range("E7").Value = "'1.00%"
$ws.Range("B8").Value = "KuCoinToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D8").Value = "'8.286"
$ws.Range("E8").Value = "'1.22%"
$ws.Range("B9").Value = "BTSEToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D9").Value = "'2.995"
$ws.Range("E9").Value = "'0.47%"
$ws.Range("B10").Value = "MXToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D10").Value = "'0.9356"
$ws.Range("E10").Value = "'2.04%"
$ws.Range("B11").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C11").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D11").Value = "'0.1082"
$ws.Range("E11").Value = "'-12.86%"
$ws.Range("B12").Value = "WazirX"
$ws.Range("C12").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D12").Value = "'0.1926"
$ws.Range("E12").Value = "'1.41%"
$ws.Range("B13").Value = "MandalaExchangeToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D13").Value = "'0.09118"
$ws.Range("E13").Value = "'2.88%"
$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D14").Value = "'0.03302"
$ws.Range("E14").Value = "'-2.43%"
$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").Value = "'0.09621"
$ws.Range("E15").Value = "'-0.88%"
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").Value = "'0.001389"
$ws.Range("E16").Value = "'1.48%"
$ws.Range("B17").Value = "TigerCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D17").Value = "'0.005882"
$ws.Range("E17").Value = "'-0.97%"
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").Value = "'3.593"
$ws.Range("E18").Value = "'1.75%"
$ws.Range("D19").Value = "'0.3410"
$ws.Range("E19").Value = "'0.04%"
$ws.Range("D20").Value = "'6.443"
$ws.Range("E20").Value = "'28.14%"
$ws.Range("D21").Value = "'0.1305"
$ws.Range("E21").Value = "'0.66%"
$ws.Range("D22").Value = "'0.2521"
$ws.Range("E22").Value = "'-2.72%"
$ws.Range("D23").Value = "'0.04414"
$ws.Range("E23").Value = "'0.57%"
$ws.Range("D24").Value = "'0.001237"
$ws.Range("E24").Value = "'1.98%"
$ws.Range("D25").Value = "'0.004624"
$ws.Range("E25").Value = "'8.88%"
$ws.Range("E26").Value = "'0.78%"
$ws.Range("D27").Value = "'0.0003992"
$ws.Range("D39").Value = "'0.02256"
$ws.Range("E39").Value = "'5.61%"
$ws.Range("D40").Value = "'0.05095"
$ws.Range("E40").Value = "'2.55%"
$ws.Range("D41").Value = "'0.007462"
$ws.Range("E41").Value = "'-3.79%"
$ws.Range("D42").Value = "'0.008884"
$ws.Range("E42").Value = "'-10.03%"
$ws.Range("E43").Value = "'0.73%"
$ws.Range("E44").Value = "'3.44%"
$ws.Range("D45").Value = "'0.009316"
$ws.Range("E45").Value = "'-3.85%"
$ws.Range("D46").Value = "'0.00006639"
$ws.Range("E46").Value = "'1.87%"
$ws.Range("E47").Value = "'0.01%"
$ws.Range("D48").Value = "'0.002862"
$ws.Range("E48").Value = "'-6.89%"
$ws.Range("D49").Value = "'0.001001"
$ws.Range("E49").Value = "'-40.75%"
$ws.Range("E50").Value = "'0.01%"
$ws.Range("E51").Value = "'0.01%"
